$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "26.639.25"
Set-TextValue "E2" "  +0.88%  "
Set-TextValue "D3" "1.641.90"
Set-TextValue "E3" "  +1.11%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "215.12"
Set-TextValue "E5" "  +1.14%  "
Set-TextValue "E6" "  +1.43%  "
Set-TextValue "E7" "  +0.00%  "
Set-TextValue "D8" "0.252"
Set-TextValue "E8" "  +1.02%  "
Set-TextValue "D9" "0.0626"
Set-TextValue "E9" "  +0.66%  "
Set-TextValue "D10" "19.12"
Set-TextValue "E10" "  +0.90%  "
Set-TextValue "D11" "0.0842"
Set-TextValue "E11" "  +0.05%  "
Set-TextValue "D12" "1.870.49"
Set-TextValue "E12" "  +1.12%  "
Set-TextValue "D13" "1.637.43"
Set-TextValue "E13" "  +1.77%  "
Set-TextValue "E14" "  +1.31%  "
Set-TextValue "E15" "  +1.55%  "
Set-TextValue "D16" "64.87"
Set-TextValue "E16" "  +0.97%  "
Set-TextValue "D17" "26.647.26"
Set-TextValue "E17" "  +0.78%  "
Set-TextValue "D18" "0.0₃0743"
Set-TextValue "E18" "  +0.67%  "
Set-TextValue "D19" "215.61"
Set-TextValue "E19" "  +0.03%  "
Set-TextValue "E20" "  -0.01%  "
Set-TextValue "D21" "4.35"
Set-TextValue "E21" "  +1.11%  "
Set-TextValue "E22" "  +0.59%  "
Set-TextValue "D23" "9.47"
Set-TextValue "E23" "  +1.96%  "
Set-TextValue "D24" "2.20"
Set-TextValue "E24" "  +12.10%  "
Set-TextValue "D25" "145.26"
Set-TextValue "E25" "  -1.67%  "
Set-TextValue "D26" "1.00"
Set-TextValue "E26" "  +0.09%  "
Set-TextValue "E27" "  +0.22%  "
Set-TextValue "E28" "  +3.90%  "
Set-TextValue "D29" "15.71"
Set-TextValue "E29" "  +0.95%  "
Set-TextValue "D30" "0.0513"
Set-TextValue "E30" "  +1.03%  "
Set-TextValue "E31" "  +1.49%  "
Set-TextValue "E32" "  +1.61%  "
Set-TextValue "D33" "3.00"
Set-TextValue "E33" "  +1.88%  "
Set-TextValue "D34" "1.277.97"
Set-TextValue "E34" "  +5.11%  "
Set-TextValue "E35" "  +3.22%  "
Set-TextValue "D36" "2.41"
Set-TextValue "E36" "  +1.48%  "
Set-TextValue "D37" "0.0178"
Set-TextValue "E37" "  +2.85%  "
Set-TextValue "D38" "0.532"
Set-TextValue "E38" "  +6.12%  "
Set-TextValue "E39" "  +3.53%  "
Set-TextValue "E40" "  +0.13%  "
Set-TextValue "D41" "0.811"
Set-TextValue "E41" "  +2.32%  "
Set-TextValue "E42" "  -0.23%  "
Set-TextValue "E43" "  +1.08%  "
Set-TextValue "D44" "1.779.95"
Set-TextValue "E44" "  +1.01%  "
Set-TextValue "D45" "91.66"
Set-TextValue "E45" "  -1.08%  "
Set-TextValue "E46" "  +8.02%  "
Set-TextValue "D47" "1.60"
Set-TextValue "E47" "  +1.13%  "
Set-TextValue "E48" "  +0.96%  "
Set-TextValue "D49" "7.73"
Set-TextValue "E49" "  +1.52%  "
Set-TextValue "D50" "0.0963"
Set-TextValue "E50" "  +1.37%  "
Set-TextValue "E51" "  -0.39%  "
